$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price-report row was inserted at row 130 (Ajo / Chino / Primera,
# Macroferia Regional de Talca), pushing all subsequent rows (130-163) down by
# one to (131-164). Insert a blank row first so the rest of the table shifts.
$ws.Rows.Item(130).Insert()

$ws.Range("A130").Value = 5
$ws.Range("B130").Value = "Macroferia Regional de Talca"
$ws.Range("C130").Value = "Maule"
$ws.Range("D130").Value = 44463
$ws.Range("E130").Value = 7
$ws.Range("F130").Value = 100112003
$ws.Range("G130").Value = "Ajo"
$ws.Range("H130").Value = "Chino"
$ws.Range("I130").Value = "Primera"
$ws.Range("J130").Value = 200
$ws.Range("K130").Value = 15000
$ws.Range("L130").Value = 15000
$ws.Range("M130").Value = 15000
$ws.Range("N130").Value = "$/malla 10 kilos"
$ws.Range("O130").Value = "China"
$ws.Range("P130").Value = 1500
$ws.Range("Q130").Value = 10
$ws.Range("R130").Value = "Hortaliza"
